$d = $word.ActiveDocument
$chi = [char]0x03C7

# 1) The header row of every "ranova" table that reports the chi-squared
#    statistic grew taller (571 -> 637 twips, i.e. 28.55pt -> 31.85pt)
#    now that the non-gaussian models were re-run with 1000 iterations.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    $headerRow = $t.Rows.Item(1)
    if ($headerRow.Range.Text.IndexOf($chi) -ge 0) {
        $headerRow.Height = 31.85
    }
}

# 2) Fix the mojibake corruption of the chi (χ) glyph used in the
#    "chi-squared" column header, everywhere it occurs in the document.
$replacement = [string]([char]0x00CF) + [string]([char]0x2021)
$null = $d.Content.Find.Execute($chi, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
